# "Adicionando mais testes" — add two new test-case rows to the
# "Tabela2" table on Planilha1 and correct the exception type recorded
# for the existing "Enviando uma String como nota" test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Tabela2")

# --- Fix the existing "Enviando uma String como nota" row --------------
# Its expected exception changes from ValueError(...) to TypeError(...).
# Row 11 is still that test case at this point (before the new row is
# inserted above it), so update it first — this also matches the order
# in which new shared strings are introduced in the saved workbook.
$ws.Cells.Item(11, 5).Value = 'TypeError("invalid note: texts are not allowed")'

# --- Insert a new row before it for "Enviando uma String ao invez de
#     uma lista" -------------------------------------------------------
$ws.Rows.Item(11).Insert()
$lo.Resize($ws.Range("A1:E13"))
$ws.Rows.Item(11).RowHeight = 18.75

$ws.Cells.Item(11, 1).Value = "escola/aluno.py"
$ws.Cells.Item(11, 2).Value = "calcular_media"
$ws.Cells.Item(11, 3).Value = "Enviando uma String ao invez de uma lista"
$ws.Cells.Item(11, 4).Value = '"olá"'
$ws.Cells.Item(11, 5).Value = 'TypeError("invalid note: texts are not allowed")'

# --- Append a new row at the end of the table for "Enviando uma nota
#     maior que 10" -----------------------------------------------------
$newRow = $lo.ListRows.Add()
$r = $newRow.Index + 1

$ws.Cells.Item($r, 1).Value = "escola/aluno.py"
$ws.Cells.Item($r, 2).Value = "calcular_media"
$ws.Cells.Item($r, 3).Value = "Enviando uma nota maior que 10"
$ws.Cells.Item($r, 4).Value = "[11.0]"
$ws.Cells.Item($r, 5).Value = 'ValueError("grades can be from 0 to 10")'

# --- Page setup / view tweaks ------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("D27").Select()
